$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse the three-line name/type/rules-text triples (rows 2-4, 5-7, 8-10,
# 11-13, 14-16) down to a single Python-tuple-style string per land, stored
# in rows 2-6.
$ws.Range("A2").Value = "('Forest', ['Basic Land — Forest', '({T}: Add {G}.)'])"
$ws.Range("A3").Value = "('Island', ['Basic Land — Island', '({T}: Add {U}.)'])"
$ws.Range("A4").Value = "('Mountain', ['Basic Land — Mountain', '({T}: Add {R}.)'])"
$ws.Range("A5").Value = "('Plains', ['Basic Land — Plains', '({T}: Add {W}.)'])"
$ws.Range("A6").Value = "('Swamp', ['Basic Land — Swamp', '({T}: Add {B}.)'])"

# The old rows 7-16 are no longer needed; delete them and shift the
# (now empty) tail up so the used range shrinks back to A1:A6.
$ws.Range("A7:A16").Delete(-4162)
